# Update simulation flux results (minimum/maximum/fluxes) after
# fixing the biomass reaction flux at 0.11 (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1856717566753563
$ws.Range("C2").Value = 0.07047413908644458
$ws.Range("D2").Value = 0.7680937323029654
$ws.Range("B3").Value = 0.2619045040918823
$ws.Range("C3").Value = -0.5570318527331276
$ws.Range("D3").Value = 0.4406610716981233
$ws.Range("B4").Value = 0.00385446058523984
$ws.Range("C4").Value = 0.003744014416556215
$ws.Range("D4").Value = 0.003856417738732076
$ws.Range("B5").Value = [double]"2.451379074012e-05"
$ws.Range("C5").Value = [double]"-3.748731746025131e-05"
$ws.Range("D5").Value = [double]"2.769017620040231e-05"
$ws.Range("B6").Value = 0.003783062672532
$ws.Range("C6").Value = 0.003778713351258491
$ws.Range("D6").Value = 0.006953723701989898
$ws.Range("B7").Value = 0.0012479245803024
$ws.Range("C7").Value = -0.001922736397396834
$ws.Range("D7").Value = 0.00192273590541154
$ws.Range("B8").Value = 0.0003374055828
$ws.Range("C8").Value = -0.0003374059085472175
$ws.Range("D8").Value = 0.0003374056625545708
$ws.Range("B9").Value = 0.007393303716672417
$ws.Range("C9").Value = 0.00738415862926447
$ws.Range("D9").Value = 0.01359176908201526
$ws.Range("B10").Value = 0.002438840229400962
$ws.Range("C10").Value = -0.003758186558964031
$ws.Range("D10").Value = 0.003758185535317971
$ws.Range("B11").Value = 0.0006594939978
$ws.Range("C11").Value = -0.0006594946534208139
$ws.Range("D11").Value = 0.0006594941536885299
$ws.Range("B12").Value = 0.003514164170316058
$ws.Range("C12").Value = 0.00351036866946521
$ws.Range("D12").Value = 0.006459437408770986
$ws.Range("B13").Value = 0.001159215330027195
$ws.Range("C13").Value = -0.001786057256931095
$ws.Range("D13").Value = 0.001786057257079671
$ws.Range("B14").Value = 0.000313420963337019
$ws.Range("C14").Value = -0.0003134209693209057
$ws.Range("D14").Value = 0.0003134209687925157
$ws.Range("B15").Value = 0.04762222545836997
$ws.Range("C15").Value = -0.04899668581452484
$ws.Range("D15").Value = 0.04899763777056636
$ws.Range("B16").Value = 0.0006874275721230985
$ws.Range("C16").Value = -0.0006870891156657752
$ws.Range("D16").Value = 0.0006888209719004614
$ws.Range("B17").Value = 0.0088082369195535
$ws.Range("C17").Value = 0.008808060179282796
$ws.Range("D17").Value = 0.01620067687608405
$ws.Range("B18").Value = 0.0029055863760462
$ws.Range("C18").Value = -0.00450011261817293
$ws.Range("D18").Value = 0.004479147180510654
$ws.Range("B19").Value = 0.00078559320015
$ws.Range("C19").Value = -0.0007905945278937255
$ws.Range("D19").Value = 0.0007867804015717369
$ws.Range("B20").Value = 0.6109670852943134
$ws.Range("C20").Value = 0.6109642550659238
$ws.Range("D20").Value = 0.6109670858853445
$ws.Range("B21").Value = 0.2493126962523645
$ws.Range("C21").Value = 0.249153365156595
$ws.Range("D21").Value = 0.3043093884402824
$ws.Range("B22").Value = 0.0024525689323982
$ws.Range("C22").Value = -0.05263003112399434
$ws.Range("D22").Value = 0.05258133794884957
$ws.Range("B23").Value = 0.003304802308192157
$ws.Range("C23").Value = -0.0495143528465315
$ws.Range("D23").Value = 0.04946565967221674
$ws.Range("B24").Value = 0.004019138646880501
$ws.Range("C24").Value = 0.004018556396839812
$ws.Range("D24").Value = 0.007387748966675581
$ws.Range("B25").Value = 0.0013257993174426
$ws.Range("C25").Value = -0.002043204413876881
$ws.Range("D25").Value = 0.002042811002352116
$ws.Range("B26").Value = 0.00035846083845
$ws.Range("C26").Value = -0.0003589442579871089
$ws.Range("D26").Value = 0.0003585508464594493
$ws.Range("B27").Value = 0.01980189811444471
$ws.Range("C27").Value = -0.07519053436097996
$ws.Range("D27").Value = 0.09287352973702369
$ws.Range("B28").Value = 0.1355565840021052
$ws.Range("C28").Value = -0.3015312766721112
$ws.Range("D28").Value = 0.673569684293486
$ws.Range("B29").Value = 0.02371014509725695
$ws.Range("C29").Value = -0.02565769299097373
$ws.Range("D29").Value = 0.02798978154973403
$ws.Range("B30").Value = 0.1644049506903602
$ws.Range("C30").Value = 0.1640875153414896
$ws.Range("D30").Value = 0.1699549874337685
$ws.Range("B31").Value = 0.04097169260652332
$ws.Range("C31").Value = -0.06400639216058052
$ws.Range("D31").Value = 0.1729259130713995
